$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change username to a numeric id, tweak the email typo, and
# bump the role from USER to MANAGER (rolecode + roletype columns).
$ws.Range("A2").Value = 123
$ws.Range("C2").Value = "123@gmail.coms"
$ws.Range("D2:E2").Value = "MANAGER"

# Row 3 (the ADMIN / user2 row) is being dropped entirely. Remove its
# mailto hyperlink first, then clear the row contents (C3 keeps its
# hyperlink-style formatting, like C4 below it, just with no value).
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -eq "mailto:testcheck@gmail.com") {
        $hl.Delete()
    }
}
$ws.Range("A3:E3").ClearContents() | Out-Null

# Update the saved selection.
$ws.Range("D8").Select() | Out-Null
